# Require angle brackets for @base and @prefix values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# Wrap the raw URL values in angle brackets to match the bracketed form
# already used elsewhere on the sheet (e.g. D5/D6/B7).
$ws.Range("D1").Value = "<http://sales.data/purchases/2015>"
$ws.Range("D2").Value = "<http://sales.data/purchases#>"
$ws.Range("D3").Value = "<http://sales.data/schema#>"

# Update the selection shown on the sheet (matches the new "after" view
# state: the C1:D3 block is highlighted with D3 as the active cell).
$ws.Activate()
$ws.Range("C1:D3").Select()
